$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column CC (81) holds the next date "1-nov", following the existing
# CB (80) column "31-oct". Copy header/format from CB so the new column
# matches the existing style, then fill in header + values.
$ws.Range("CB1:CB11").Copy()
$ws.Range("CC1:CC11").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("CC1").Value = "1-nov"

$values = @(7, 10, 8, 3, 11, 5, 13, 11, 6, 0)
for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 81).Value = $values[$i]
}

$ws.Range("CC10").Select() | Out-Null
